$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Header row (row 1) — add the new "capacity" column (C) and the
# standard metadata columns (H:N) that every other sheet already has.
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Data row (row 2) — fill in the metadata columns for the car record.
$ws.Cells.Item(2, 8).Value = "land"
$ws.Cells.Item(2, 9).Value = "normal"
$ws.Cells.Item(2, 10).Value = "2012-04-30"
$ws.Cells.Item(2, 11).Value = "劉櫂豪"
$ws.Cells.Item(2, 12).Value = 1762
$ws.Cells.Item(2, 13).Value = "tmpba991"
$ws.Cells.Item(2, 14).Value = 29

# Apply the header/body styles used elsewhere in the workbook to the
# newly populated cells so they match the existing columns.
$ws.Range("C1:N1").Style = $ws.Range("B1").Style
$ws.Range("H2:N2").Style = $ws.Range("B2").Style
